$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 262.4375
$ws.Range("J55").Value = 225.875
$ws.Range("L55").Value = 225.875
$ws.Range("N55").Value = -653.875

# Row 107
$ws.Range("H107").Value = 3970.2563
$ws.Range("I107").Value = 4734.129
$ws.Range("J107").Value = 1010.25
$ws.Range("K107").Value = 4734.129
$ws.Range("L107").Value = 1010.25
$ws.Range("M107").Value = -2814.129
$ws.Range("N107").Value = -4850.25

# Row 132
$ws.Range("H132").Value = 4242861.5
$ws.Range("I132").Value = 5057.6665
$ws.Range("J132").Value = 10875946
$ws.Range("K132").Value = 15172.9995
$ws.Range("L132").Value = 32627838
$ws.Range("M132").Value = -12642.9995
$ws.Range("N132").Value = -32632898

# Row 138
$ws.Range("H138").Value = 3970280.8
$ws.Range("I138").Value = 1759.9445
$ws.Range("J138").Value = 9261642
$ws.Range("K138").Value = 5279.833500000001
$ws.Range("L138").Value = 27784926
$ws.Range("M138").Value = -139.8335000000006
$ws.Range("N138").Value = -27795206

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8788.016
$ws.Range("I32").Value = 9074.709000000001
$ws.Range("K32").Value = 9074.709000000001
$ws.Range("M32").Value = -8787.709000000001

# Row 97
$ws.Range("H97").Value = 11405.1
$ws.Range("I97").Value = 18415
$ws.Range("J97").Value = 890.25
$ws.Range("K97").Value = 18415
$ws.Range("L97").Value = 890.25
$ws.Range("M97").Value = -17919
$ws.Range("N97").Value = -1882.25

# Row 102
$ws.Range("H102").Value = 12110
$ws.Range("I102").Value = 11600
$ws.Range("J102").Value = 12747.5
$ws.Range("K102").Value = 11600
$ws.Range("L102").Value = 12747.5
$ws.Range("M102").Value = -9978
$ws.Range("N102").Value = -15991.5

# Row 123
$ws.Range("H123").Value = 33429
$ws.Range("J123").Value = 33429
$ws.Range("L123").Value = 33429
$ws.Range("N123").Value = -43229

$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 15999.857
$ws.Range("I97").Value = 10333.167
$ws.Range("J97").Value = 50000
$ws.Range("K97").Value = 10333.167
$ws.Range("L97").Value = 50000
$ws.Range("M97").Value = -9342.166999999999
$ws.Range("N97").Value = -51982

# Row 99
$ws.Range("H99").Value = 2094
$ws.Range("I99").Value = 2034.1666
$ws.Range("K99").Value = 2034.1666
$ws.Range("M99").Value = -536.1666

# Row 103
$ws.Range("H103").Value = 27894.5
$ws.Range("J103").Value = 27894.5
$ws.Range("L103").Value = 27894.5
$ws.Range("N103").Value = -30238.5

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 2273.4546
$ws.Range("I86").Value = 2460
$ws.Range("J86").Value = 2118
$ws.Range("K86").Value = 2460
$ws.Range("L86").Value = 2118
$ws.Range("M86").Value = -1337
$ws.Range("N86").Value = -4364

# Row 89
$ws.Range("H89").Value = 2273.4546
$ws.Range("I89").Value = 2460
$ws.Range("J89").Value = 2118
$ws.Range("K89").Value = 12300
$ws.Range("L89").Value = 10590
$ws.Range("M89").Value = -6684
$ws.Range("N89").Value = -21822

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 877.6941
$ws.Range("I68").Value = 601.0492
$ws.Range("J68").Value = 1580.8334
$ws.Range("K68").Value = 1803.1476
$ws.Range("L68").Value = 4742.5002
$ws.Range("M68").Value = -992.1476000000002
$ws.Range("N68").Value = -6364.5002

# Row 71
$ws.Range("H71").Value = 877.6941
$ws.Range("I71").Value = 601.0492
$ws.Range("J71").Value = 1580.8334
$ws.Range("K71").Value = 5409.442800000001
$ws.Range("L71").Value = 14227.5006
$ws.Range("M71").Value = -1353.442800000001
$ws.Range("N71").Value = -22339.5006

# Row 113
$ws.Range("H113").Value = 1193.5385
$ws.Range("I113").Value = 839.65216
$ws.Range("J113").Value = 1702.25
$ws.Range("K113").Value = 2518.95648
$ws.Range("L113").Value = 5106.75
$ws.Range("M113").Value = -348.9564799999998
$ws.Range("N113").Value = -9446.75

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 39963.707
$ws.Range("I70").Value = 89852.28999999999
$ws.Range("J70").Value = 5041.7
$ws.Range("K70").Value = 89852.28999999999
$ws.Range("L70").Value = 5041.7
$ws.Range("M70").Value = -89582.28999999999
$ws.Range("N70").Value = -5581.7

# Row 73
$ws.Range("H73").Value = 39963.707
$ws.Range("I73").Value = 89852.28999999999
$ws.Range("J73").Value = 5041.7
$ws.Range("K73").Value = 89852.28999999999
$ws.Range("L73").Value = 5041.7
$ws.Range("M73").Value = -88916.28999999999
$ws.Range("N73").Value = -6913.7

# Row 97
$ws.Range("H97").Value = 766.6667
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -4
$ws.Range("N97").Value = -1892

# Row 132
$ws.Range("H132").Value = 5394.2915
$ws.Range("I132").Value = 4094.3635
$ws.Range("J132").Value = 6494.231
$ws.Range("K132").Value = 12283.0905
$ws.Range("L132").Value = 19482.693
$ws.Range("M132").Value = -9753.0905
$ws.Range("N132").Value = -24542.693

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 961.1818
$ws.Range("I46").Value = 744.4286
$ws.Range("J46").Value = 1340.5
$ws.Range("K46").Value = 744.4286
$ws.Range("L46").Value = 1340.5
$ws.Range("M46").Value = -556.4286
$ws.Range("N46").Value = -1716.5

# Row 55
$ws.Range("H55").Value = 269.12
$ws.Range("I55").Value = 120.545456
$ws.Range("J55").Value = 385.85715
$ws.Range("K55").Value = 120.545456
$ws.Range("L55").Value = 385.85715
$ws.Range("M55").Value = 52.454544
$ws.Range("N55").Value = -731.85715

# Row 100
$ws.Range("H100").Value = 2040.8
$ws.Range("I100").Value = 1900
$ws.Range("K100").Value = 1900
$ws.Range("M100").Value = -1359

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 11831
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 12891.182
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 12891.182
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -14139.182

# Row 65
$ws.Range("H65").Value = 11831
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 12891.182
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 64455.91
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -70695.91

# Row 74
$ws.Range("H74").Value = 10687.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10687.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 10687.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -12559.5

# Row 77
$ws.Range("H77").Value = 10687.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10687.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 32062.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -41422.5

# Row 96
$ws.Range("H96").Value = 4271.4287
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 4483.3335
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 4483.3335
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -7229.3335

# Row 123
$ws.Range("H123").Value = 36883.5
$ws.Range("J123").Value = 36883.5
$ws.Range("L123").Value = 36883.5
$ws.Range("N123").Value = -46683.5
